$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (Volume/Number + date range), rich-text runs collapse to plain text
# (all runs share identical formatting, so no visual change results) ---
$ws.Cells.Item(8, 1).Value = "Volume 31   Number  50"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# --- Crime complaint statistics table updates (rows 14-33) ---

# Row 14
$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "0"

# Row 15
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 9).Value = 18
$ws.Cells.Item(15, 11).Value = 125
$ws.Cells.Item(15, 12).Value = 125
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = 0

# Row 16
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 7).Value = 3
$ws.Cells.Item(16, 8).Value = -66.666666666666
$ws.Cells.Item(16, 9).Value = 52
$ws.Cells.Item(16, 10).Value = 48
$ws.Cells.Item(16, 11).Value = 8.333333333333
$ws.Cells.Item(16, 12).Value = -16.129032258064
$ws.Cells.Item(16, 13).Value = -52.727272727272
$ws.Cells.Item(16, 14).Value = -82.608695652173

# Row 17
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 4
$ws.Cells.Item(17, 5).Value = -75
$ws.Cells.Item(17, 6).Value = 17
$ws.Cells.Item(17, 7).Value = 11
$ws.Cells.Item(17, 8).Value = 54.545454545454
$ws.Cells.Item(17, 9).Value = 159
$ws.Cells.Item(17, 10).Value = 163
$ws.Cells.Item(17, 11).Value = -2.453987730061
$ws.Cells.Item(17, 12).Value = 30.327868852459
$ws.Cells.Item(17, 13).Value = 19.548872180451
$ws.Cells.Item(17, 14).Value = -45.547945205479

# Row 18
$ws.Cells.Item(18, 4).Value = 4
$ws.Cells.Item(18, 10).Value = 97
$ws.Cells.Item(18, 11).Value = -3.092783505154
$ws.Cells.Item(18, 12).Value = 13.253012048192
$ws.Cells.Item(18, 13).Value = -52.763819095477
$ws.Cells.Item(18, 14).Value = -92.780337941628

# Row 19
$ws.Cells.Item(19, 4).Value = 10
$ws.Cells.Item(19, 5).Value = -20
$ws.Cells.Item(19, 6).Value = 31
$ws.Cells.Item(19, 7).Value = 26
$ws.Cells.Item(19, 8).Value = 19.230769230769
$ws.Cells.Item(19, 9).Value = 445
$ws.Cells.Item(19, 10).Value = 456
$ws.Cells.Item(19, 11).Value = -2.412280701754
$ws.Cells.Item(19, 12).Value = 24.649859943977
$ws.Cells.Item(19, 13).Value = 13.520408163265
$ws.Cells.Item(19, 14).Value = -44.993819530284

# Row 20
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 100
$ws.Cells.Item(20, 6).Value = 4
$ws.Cells.Item(20, 8).Value = -42.857142857142
$ws.Cells.Item(20, 9).Value = 75
$ws.Cells.Item(20, 10).Value = 106
$ws.Cells.Item(20, 11).Value = -29.245283018867
$ws.Cells.Item(20, 12).Value = -43.181818181818
$ws.Cells.Item(20, 13).Value = -33.62831858407
$ws.Cells.Item(20, 14).Value = -97.2087830294

# Row 21
$ws.Cells.Item(21, 3).Value = 13
$ws.Cells.Item(21, 4).Value = 20
$ws.Cells.Item(21, 5).Value = -35
$ws.Cells.Item(21, 6).Value = 58
$ws.Cells.Item(21, 8).Value = -7.936507936507
$ws.Cells.Item(21, 9).Value = 844
$ws.Cells.Item(21, 10).Value = 880
$ws.Cells.Item(21, 11).Value = -4.090909090909
$ws.Cells.Item(21, 12).Value = 10.32679738562
$ws.Cells.Item(21, 13).Value = -12.71975180972
$ws.Cells.Item(21, 14).Value = -84.399260628465

# Row 23
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0"
$ws.Cells.Item(23, 5).Value = "***.*"
$ws.Cells.Item(23, 6).Value = 4
$ws.Cells.Item(23, 7).Value = 1
$ws.Cells.Item(23, 8).Value = 300
$ws.Cells.Item(23, 9).Value = 27
$ws.Cells.Item(23, 11).Value = -18.181818181818
$ws.Cells.Item(23, 12).Value = 35
$ws.Cells.Item(23, 13).Value = 17.391304347826

# Row 24
$ws.Cells.Item(24, 3).Value = 11
$ws.Cells.Item(24, 4).Value = 21
$ws.Cells.Item(24, 5).Value = -47.619047619047
$ws.Cells.Item(24, 6).Value = 52
$ws.Cells.Item(24, 7).Value = 99
$ws.Cells.Item(24, 8).Value = -47.474747474747
$ws.Cells.Item(24, 9).Value = 947
$ws.Cells.Item(24, 10).Value = 1056
$ws.Cells.Item(24, 11).Value = -10.321969696969
$ws.Cells.Item(24, 12).Value = 13.822115384615
$ws.Cells.Item(24, 13).Value = -42.150274893097

# Row 25
$ws.Cells.Item(25, 3).Value = 4
$ws.Cells.Item(25, 5).Value = -33.333333333333
$ws.Cells.Item(25, 6).Value = 29
$ws.Cells.Item(25, 7).Value = 22
$ws.Cells.Item(25, 8).Value = 31.818181818181
$ws.Cells.Item(25, 9).Value = 465
$ws.Cells.Item(25, 10).Value = 418
$ws.Cells.Item(25, 11).Value = 11.244019138756
$ws.Cells.Item(25, 12).Value = 113.302752293578

# Row 26
$ws.Cells.Item(26, 3).Value = 12
$ws.Cells.Item(26, 4).Value = 7
$ws.Cells.Item(26, 5).Value = 71.428571428571
$ws.Cells.Item(26, 6).Value = 28
$ws.Cells.Item(26, 7).Value = 26
$ws.Cells.Item(26, 8).Value = 7.692307692307
$ws.Cells.Item(26, 9).Value = 350
$ws.Cells.Item(26, 10).Value = 315
$ws.Cells.Item(26, 11).Value = 11.111111111111
$ws.Cells.Item(26, 12).Value = 8.024691358024
$ws.Cells.Item(26, 13).Value = -32.562620423892

# Row 27
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 9).Value = 23
$ws.Cells.Item(27, 11).Value = 64.285714285714
$ws.Cells.Item(27, 12).Value = 21.052631578947

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0"
$ws.Cells.Item(28, 5).Value = "***.*"
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(28, 12).Value = 6.451612903225

# Row 31
$ws.Cells.Item(31, 12).Value = -25

# Row 33
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 9).Value = 5
$ws.Cells.Item(33, 11).Value = -28.571428571428
$ws.Cells.Item(33, 12).Value = 150
